# update lowongan kerja user
# - Row 2 (NINDYA RIZQY): birth-place/date and address text corrected
# - Row 3 (IDA FANIA): gender, email, birth date/place, address, postal code
#   and status filled in / corrected
# - Row 4 (LIA ETIKASARI): email corrected (gender text re-used from shared pool)
# - Row 5 (KURNIA AINUN): gender, email, birth date/place, address and
#   postal code filled in / corrected

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : NINDYA RIZQY -------------------------------------------------
$ws.Range("E2").Value = "temanggung abc, 2000-02-18"
$ws.Range("F2").Value = "KENANGAN MANTAN abcDE, SINUNUKAN, KABUPATEN MANDAILING NATAL, SUMATERA UTARA"

# --- Row 3 : IDA FANIA -----------------------------------------------------
$ws.Range("B3").Value = "Perempuan"
$ws.Range("C3").Value = "idafania33@gmail.com"
$ws.Range("E3").Value = "Jakarta, 2025-11-06"
$ws.Range("F3").Value = "Pungangan, MAJENANG, KABUPATEN CILACAP, JAWA TENGAH"
$ws.Range("G3").Value = 56351
$ws.Range("I3").Value = "AKTIF"

# --- Row 4 : LIA ETIKASARI ---------------------------------------------------
$ws.Range("B4").Value = "Tidak Ingin Menyebutkan"
$ws.Range("C4").Value = "liaetikasari0826@gmail.com"

# --- Row 5 : KURNIA AINUN ---------------------------------------------------
$ws.Range("B5").Value = "Transgender"
$ws.Range("C5").Value = "etikasarilia26@gmail.com"
$ws.Range("E5").Value = "Temanggung, 2003-03-03"
$ws.Range("F5").Value = "KENANGAN, KOTAGEDE, KOTA YOGYAKARTA, DI YOGYAKARTA"
$ws.Range("G5").Value = 12345

# --- Column widths (Alamat/Tanggal lahir columns grew wider with new text) -
# ColumnWidth setter pads by 5/6 (0.8333) in the saved XML, so subtract that
# back out to land on the exact target widths of 31 / 90.
$ws.Columns.Item(5).ColumnWidth = 31 - 5/6
$ws.Columns.Item(6).ColumnWidth = 90 - 5/6
